$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 322360.0138707894
$ws.Range("D2").Value = 0.2634730538922156
$ws.Range("E2").Value = 0.2645290581162325
$ws.Range("F2").Value = 0.264
$ws.Range("G2").Value = 0.0008205980537835252

$ws.Range("C3").Value = 51391070.5775004
$ws.Range("D3").Value = 0.4445336544894955
$ws.Range("E3").Value = 0.4428152492668622
$ws.Range("F3").Value = 0.4436727879799666
$ws.Range("G3").Value = 0.0000645163976409918

$ws.Range("C4").Value = 25881672.87276903
$ws.Range("D4").Value = 0.4438381687516635
$ws.Range("E4").Value = 0.4445481205011997
$ws.Range("F4").Value = 0.4441928609483218
$ws.Range("G4").Value = 0.000128958633144512

$ws.Range("C5").Value = 51651234.30294618
$ws.Range("D5").Value = 0.8387912673056444
$ws.Range("E5").Value = 0.8399093575046654
$ws.Range("F5").Value = 0.8393499400559479
$ws.Range("G5").Value = 0.0001220725280778166

$ws.Range("C6").Value = 40669995.20361029
$ws.Range("D6").Value = 0.4445336544894955
$ws.Range("E6").Value = 0.4428152492668622
$ws.Range("F6").Value = 0.4436727879799666
$ws.Range("G6").Value = 0.00008152365713286255

$ws.Range("C7").Value = 25873986.65569574
$ws.Range("D7").Value = 0.8386366662228731
$ws.Range("E7").Value = 0.839642761930152
$ws.Range("F7").Value = 0.8391394125091587
$ws.Range("G7").Value = 0.00024359509074605

$ws.Range("C8").Value = 40875827.43311323
$ws.Range("D8").Value = 0.8387912673056444
$ws.Range("E8").Value = 0.8399093575046654
$ws.Range("F8").Value = 0.8393499400559479
$ws.Range("G8").Value = 0.0001542524554400209

$ws.Range("C9").Value = 20487884.13131923
$ws.Range("D9").Value = 0.4438381687516635
$ws.Range("E9").Value = 0.4445481205011997
$ws.Range("F9").Value = 0.4441928609483218
$ws.Range("G9").Value = 0.0001629092167728291

$ws.Range("C10").Value = 20481803.0679683
$ws.Range("D10").Value = 0.8386366662228731
$ws.Range("E10").Value = 0.839642761930152
$ws.Range("F10").Value = 0.8391394125091587
$ws.Range("G10").Value = 0.0003077256482957434

